$d = $word.ActiveDocument

$d.Content.Find.Execute("71+8=79", $true, $false, $false, $false, $false, $true, 1, $false, "60-35=25", 2) | Out-Null
$d.Content.Find.Execute("57-6=51", $true, $false, $false, $false, $false, $true, 1, $false, "28-6=22", 2) | Out-Null
$d.Content.Find.Execute("82+6=88", $true, $false, $false, $false, $false, $true, 1, $false, "70-12=58", 2) | Out-Null
$d.Content.Find.Execute("74+7=81", $true, $false, $false, $false, $false, $true, 1, $false, "56-9=47", 2) | Out-Null
$d.Content.Find.Execute("98-23=75", $true, $false, $false, $false, $false, $true, 1, $false, "16-10=6", 2) | Out-Null
$d.Content.Find.Execute("88-41=47", $true, $false, $false, $false, $false, $true, 1, $false, "41-27=14", 2) | Out-Null
$d.Content.Find.Execute("4+1=5", $true, $false, $false, $false, $false, $true, 1, $false, "10+64=74", 2) | Out-Null
$d.Content.Find.Execute("81-46=35", $true, $false, $false, $false, $false, $true, 1, $false, "71-49=22", 2) | Out-Null
$d.Content.Find.Execute("92-75=17", $true, $false, $false, $false, $false, $true, 1, $false, "77+20=97", 2) | Out-Null
$d.Content.Find.Execute("1+52=53", $true, $false, $false, $false, $false, $true, 1, $false, "2+72=74", 2) | Out-Null
$d.Content.Find.Execute("97+1=98", $true, $false, $false, $false, $false, $true, 1, $false, "31+63=94", 2) | Out-Null
$d.Content.Find.Execute("1+50=51", $true, $false, $false, $false, $false, $true, 1, $false, "87-83=4", 2) | Out-Null
$d.Content.Find.Execute("22+23=45", $true, $false, $false, $false, $false, $true, 1, $false, "51-3=48", 2) | Out-Null
$d.Content.Find.Execute("78+3=81", $true, $false, $false, $false, $false, $true, 1, $false, "43-31=12", 2) | Out-Null
$d.Content.Find.Execute("15+84=99", $true, $false, $false, $false, $false, $true, 1, $false, "41-3=38", 2) | Out-Null
$d.Content.Find.Execute("69-26=43", $true, $false, $false, $false, $false, $true, 1, $false, "57+33=90", 2) | Out-Null
$d.Content.Find.Execute("84-56=28", $true, $false, $false, $false, $false, $true, 1, $false, "76-7=69", 2) | Out-Null
$d.Content.Find.Execute("74+23=97", $true, $false, $false, $false, $false, $true, 1, $false, "86-44=42", 2) | Out-Null
$d.Content.Find.Execute("11+34=45", $true, $false, $false, $false, $false, $true, 1, $false, "52-23=29", 2) | Out-Null
$d.Content.Find.Execute("12+85=97", $true, $false, $false, $false, $false, $true, 1, $false, "88-44=44", 2) | Out-Null
$d.Content.Find.Execute("21+56=77", $true, $false, $false, $false, $false, $true, 1, $false, "17-15=2", 2) | Out-Null
$d.Content.Find.Execute("20+72=92", $true, $false, $false, $false, $false, $true, 1, $false, "24+11=35", 2) | Out-Null
$d.Content.Find.Execute("64+2=66", $true, $false, $false, $false, $false, $true, 1, $false, "6+28=34", 2) | Out-Null
$d.Content.Find.Execute("16+73=89", $true, $false, $false, $false, $false, $true, 1, $false, "39-31=8", 2) | Out-Null
$d.Content.Find.Execute("87-29=58", $true, $false, $false, $false, $false, $true, 1, $false, "80+15=95", 2) | Out-Null
$d.Content.Find.Execute("79-72=7", $true, $false, $false, $false, $false, $true, 1, $false, "48-36=12", 2) | Out-Null
$d.Content.Find.Execute("46-24=22", $true, $false, $false, $false, $false, $true, 1, $false, "63-0=63", 2) | Out-Null
$d.Content.Find.Execute("58-0=58", $true, $false, $false, $false, $false, $true, 1, $false, "83-23=60", 2) | Out-Null
$d.Content.Find.Execute("93-30=63", $true, $false, $false, $false, $false, $true, 1, $false, "84-25=59", 2) | Out-Null
$d.Content.Find.Execute("18+16=34", $true, $false, $false, $false, $false, $true, 1, $false, "12+7=19", 2) | Out-Null
$d.Content.Find.Execute("15+79=94", $true, $false, $false, $false, $false, $true, 1, $false, "15+52=67", 2) | Out-Null
$d.Content.Find.Execute("72+25=97", $true, $false, $false, $false, $false, $true, 1, $false, "54-27=27", 2) | Out-Null
$d.Content.Find.Execute("3+62=65", $true, $false, $false, $false, $false, $true, 1, $false, "80-29=51", 2) | Out-Null
$d.Content.Find.Execute("63+33=96", $true, $false, $false, $false, $false, $true, 1, $false, "47+51=98", 2) | Out-Null
$d.Content.Find.Execute("30+36=66", $true, $false, $false, $false, $false, $true, 1, $false, "4+80=84", 2) | Out-Null
$d.Content.Find.Execute("89-58=31", $true, $false, $false, $false, $false, $true, 1, $false, "0+42=42", 2) | Out-Null
$d.Content.Find.Execute("65+16=81", $true, $false, $false, $false, $false, $true, 1, $false, "18-10=8", 2) | Out-Null
$d.Content.Find.Execute("76-68=8", $true, $false, $false, $false, $false, $true, 1, $false, "75+5=80", 2) | Out-Null
$d.Content.Find.Execute("23+51=74", $true, $false, $false, $false, $false, $true, 1, $false, "52+13=65", 2) | Out-Null
$d.Content.Find.Execute("26+6=32", $true, $false, $false, $false, $false, $true, 1, $false, "96-60=36", 2) | Out-Null
$d.Content.Find.Execute("94-57=37", $true, $false, $false, $false, $false, $true, 1, $false, "63-59=4", 2) | Out-Null
$d.Content.Find.Execute("94-43=51", $true, $false, $false, $false, $false, $true, 1, $false, "31+64=95", 2) | Out-Null
$d.Content.Find.Execute("28-20=8", $true, $false, $false, $false, $false, $true, 1, $false, "86-67=19", 2) | Out-Null
$d.Content.Find.Execute("35+1=36", $true, $false, $false, $false, $false, $true, 1, $false, "60-31=29", 2) | Out-Null
$d.Content.Find.Execute("50+32=82", $true, $false, $false, $false, $false, $true, 1, $false, "0+69=69", 2) | Out-Null
$d.Content.Find.Execute("3+71=74", $true, $false, $false, $false, $false, $true, 1, $false, "47-20=27", 2) | Out-Null
$d.Content.Find.Execute("42-42=0", $true, $false, $false, $false, $false, $true, 1, $false, "95-73=22", 2) | Out-Null
$d.Content.Find.Execute("76+7=83", $true, $false, $false, $false, $false, $true, 1, $false, "45+7=52", 2) | Out-Null
$d.Content.Find.Execute("29+68=97", $true, $false, $false, $false, $false, $true, 1, $false, "78+20=98", 2) | Out-Null
$d.Content.Find.Execute("58-1=57", $true, $false, $false, $false, $false, $true, 1, $false, "51-44=7", 2) | Out-Null
$d.Content.Find.Execute("33+65=98", $true, $false, $false, $false, $false, $true, 1, $false, "99-45=54", 2) | Out-Null
$d.Content.Find.Execute("37-26=11", $true, $false, $false, $false, $false, $true, 1, $false, "51-1=50", 2) | Out-Null
$d.Content.Find.Execute("81-58=23", $true, $false, $false, $false, $false, $true, 1, $false, "31+25=56", 2) | Out-Null
$d.Content.Find.Execute("20+55=75", $true, $false, $false, $false, $false, $true, 1, $false, "33+37=70", 2) | Out-Null
$d.Content.Find.Execute("85-83=2", $true, $false, $false, $false, $false, $true, 1, $false, "36+36=72", 2) | Out-Null
$d.Content.Find.Execute("83-22=61", $true, $false, $false, $false, $false, $true, 1, $false, "18+45=63", 2) | Out-Null
$d.Content.Find.Execute("23-11=12", $true, $false, $false, $false, $false, $true, 1, $false, "53-19=34", 2) | Out-Null
$d.Content.Find.Execute("83+1=84", $true, $false, $false, $false, $false, $true, 1, $false, "86+13=99", 2) | Out-Null
$d.Content.Find.Execute("54+20=74", $true, $false, $false, $false, $false, $true, 1, $false, "57+33=90", 2) | Out-Null
$d.Content.Find.Execute("32+65=97", $true, $false, $false, $false, $false, $true, 1, $false, "9+54=63", 2) | Out-Null
$d.Content.Find.Execute("7+26=33", $true, $false, $false, $false, $false, $true, 1, $false, "72+26=98", 2) | Out-Null
$d.Content.Find.Execute("96-85=11", $true, $false, $false, $false, $false, $true, 1, $false, "57-30=27", 2) | Out-Null
$d.Content.Find.Execute("17+0=17", $true, $false, $false, $false, $false, $true, 1, $false, "3+10=13", 2) | Out-Null
$d.Content.Find.Execute("42+5=47", $true, $false, $false, $false, $false, $true, 1, $false, "78-58=20", 2) | Out-Null
$d.Content.Find.Execute("47+5=52", $true, $false, $false, $false, $false, $true, 1, $false, "61-50=11", 2) | Out-Null
$d.Content.Find.Execute("90-27=63", $true, $false, $false, $false, $false, $true, 1, $false, "99-67=32", 2) | Out-Null
$d.Content.Find.Execute("81-75=6", $true, $false, $false, $false, $false, $true, 1, $false, "63-61=2", 2) | Out-Null
$d.Content.Find.Execute("55+33=88", $true, $false, $false, $false, $false, $true, 1, $false, "85-30=55", 2) | Out-Null
$d.Content.Find.Execute("48+12=60", $true, $false, $false, $false, $false, $true, 1, $false, "87+11=98", 2) | Out-Null
$d.Content.Find.Execute("79-56=23", $true, $false, $false, $false, $false, $true, 1, $false, "30+32=62", 2) | Out-Null
$d.Content.Find.Execute("97-57=40", $true, $false, $false, $false, $false, $true, 1, $false, "56-32=24", 2) | Out-Null
$d.Content.Find.Execute("36-6=30", $true, $false, $false, $false, $false, $true, 1, $false, "17-1=16", 2) | Out-Null
$d.Content.Find.Execute("66-46=20", $true, $false, $false, $false, $false, $true, 1, $false, "52+27=79", 2) | Out-Null
$d.Content.Find.Execute("69-23=46", $true, $false, $false, $false, $false, $true, 1, $false, "84-61=23", 2) | Out-Null
$d.Content.Find.Execute("79+12=91", $true, $false, $false, $false, $false, $true, 1, $false, "82-23=59", 2) | Out-Null
$d.Content.Find.Execute("12+70=82", $true, $false, $false, $false, $false, $true, 1, $false, "83-37=46", 2) | Out-Null
$d.Content.Find.Execute("94-9=85", $true, $false, $false, $false, $false, $true, 1, $false, "14+33=47", 2) | Out-Null
$d.Content.Find.Execute("73-39=34", $true, $false, $false, $false, $false, $true, 1, $false, "91-49=42", 2) | Out-Null
$d.Content.Find.Execute("98-65=33", $true, $false, $false, $false, $false, $true, 1, $false, "47-13=34", 2) | Out-Null
$d.Content.Find.Execute("32+11=43", $true, $false, $false, $false, $false, $true, 1, $false, "82-64=18", 2) | Out-Null
$d.Content.Find.Execute("17+37=54", $true, $false, $false, $false, $false, $true, 1, $false, "10+6=16", 2) | Out-Null
$d.Content.Find.Execute("73-56=17", $true, $false, $false, $false, $false, $true, 1, $false, "26+4=30", 2) | Out-Null
$d.Content.Find.Execute("94-76=18", $true, $false, $false, $false, $false, $true, 1, $false, "63-44=19", 2) | Out-Null
$d.Content.Find.Execute("77-22=55", $true, $false, $false, $false, $false, $true, 1, $false, "8+62=70", 2) | Out-Null
$d.Content.Find.Execute("6+11=17", $true, $false, $false, $false, $false, $true, 1, $false, "15+61=76", 2) | Out-Null
$d.Content.Find.Execute("67-3=64", $true, $false, $false, $false, $false, $true, 1, $false, "42-7=35", 2) | Out-Null
$d.Content.Find.Execute("23+28=51", $true, $false, $false, $false, $false, $true, 1, $false, "51+43=94", 2) | Out-Null
$d.Content.Find.Execute("60-26=34", $true, $false, $false, $false, $false, $true, 1, $false, "94-25=69", 2) | Out-Null
$d.Content.Find.Execute("18+29=47", $true, $false, $false, $false, $false, $true, 1, $false, "43+41=84", 2) | Out-Null
$d.Content.Find.Execute("82-81=1", $true, $false, $false, $false, $false, $true, 1, $false, "50-1=49", 2) | Out-Null
$d.Content.Find.Execute("97-52=45", $true, $false, $false, $false, $false, $true, 1, $false, "48-20=28", 2) | Out-Null
$d.Content.Find.Execute("62-52=10", $true, $false, $false, $false, $false, $true, 1, $false, "57+8=65", 2) | Out-Null
$d.Content.Find.Execute("15-4=11", $true, $false, $false, $false, $false, $true, 1, $false, "11+27=38", 2) | Out-Null
$d.Content.Find.Execute("69+22=91", $true, $false, $false, $false, $false, $true, 1, $false, "26-26=0", 2) | Out-Null
$d.Content.Find.Execute("85-20=65", $true, $false, $false, $false, $false, $true, 1, $false, "97-74=23", 2) | Out-Null
$d.Content.Find.Execute("98-9=89", $true, $false, $false, $false, $false, $true, 1, $false, "74+25=99", 2) | Out-Null
$d.Content.Find.Execute("83-75=8", $true, $false, $false, $false, $false, $true, 1, $false, "45-14=31", 2) | Out-Null
$d.Content.Find.Execute("42+56=98", $true, $false, $false, $false, $false, $true, 1, $false, "88-40=48", 2) | Out-Null
$d.Content.Find.Execute("35-15=20", $true, $false, $false, $false, $false, $true, 1, $false, "74+2=76", 2) | Out-Null
$d.Content.Find.Execute("87-6=81", $true, $false, $false, $false, $false, $true, 1, $false, "93-51=42", 2) | Out-Null
